$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Theme color swap: accent1 <-> accent5 ---
# Office theme color indices (ThemeColorScheme.Colors): 1 dk1, 2 lt1, 3 dk2, 4 lt2,
# 5 accent1, 6 accent2, 7 accent3, 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$cs = $wb.Theme.ThemeColorScheme
$accent1 = $cs.Colors(5)
$accent5 = $cs.Colors(9)
$accent1.RGB = 12874308   # BGR for 4472C4
$accent5.RGB = 13998939   # BGR for 5B9BD5

# --- Remove the extra sheets, keep only Sheet1 ---
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Sheet2").Delete()
$wb.Worksheets.Item("Sheet3").Delete()
$excel.DisplayAlerts = $true

# --- Update the data values in B2:D4 ---
$ws.Range("B2").Value = 0.88061396996660235
$ws.Range("C2").Value = 0.42222484900448654
$ws.Range("D2").Value = -0.21504700133412391

$ws.Range("B3").Value = 0.44624828865828775
$ws.Range("C3").Value = -0.58641965468809454
$ws.Range("D3").Value = 0.67599885611223154

$ws.Range("B4").Value = 0.15931572668512545
$ws.Range("C4").Value = -0.69125839270032352
$ws.Range("D4").Value = -0.704826456478583

# --- Select the used range A1:D4 on Sheet1 to match the saved selection ---
$ws.Activate()
$ws.Range("A1:D4").Select()
